$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)
$rng = $ws.Range("A33:C33")
Write-Host $rng.Borders.Item(7).LineStyle
